$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab name does not need to match, but keep consistent)
$ws.Name = "R4RResourceLoad"

# Update the C2/C3 values: strip the leading ": " prefix from the descriptive text
$ws.Range("C2").Value = "Enhancer Linking by Methylation/Expression Relationships (ELMER) - National Cancer Institute"
$ws.Range("C3").Value = "Next-Generation Clustered Heat Maps (NG-CHM) - National Cancer Institute"
